$wb = $excel.ActiveWorkbook

# New row 50 data for each of the 4 sheets (A: date+style, B-E: inline hex strings, F-I: numbers)
$rowDate = 45836.43600694444

$sheet1 = $wb.Worksheets.Item("DE_LFT_#1")
$sheet1.Cells.Item(50, 1).Value = $rowDate
$sheet1.Cells.Item(50, 1).NumberFormat = $sheet1.Cells.Item(49, 1).NumberFormat
$sheet1.Cells.Item(50, 2).Value = "0x01,0x7c"
$sheet1.Cells.Item(50, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$sheet1.Cells.Item(50, 4).Value = "0x01,0x64"
$sheet1.Cells.Item(50, 5).Value = "0x14"
$sheet1.Cells.Item(50, 6).Value = 380
$sheet1.Cells.Item(50, 7).Value = [double]"7.598631275147109e+23"
$sheet1.Cells.Item(50, 8).Value = 356
$sheet1.Cells.Item(50, 9).Value = 14

$sheet2 = $wb.Worksheets.Item("DE_LFT_#2")
$sheet2.Cells.Item(50, 1).Value = $rowDate
$sheet2.Cells.Item(50, 1).NumberFormat = $sheet2.Cells.Item(49, 1).NumberFormat
$sheet2.Cells.Item(50, 2).Value = "0x01,0x7c"
$sheet2.Cells.Item(50, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$sheet2.Cells.Item(50, 4).Value = "0x01,0x64"
$sheet2.Cells.Item(50, 5).Value = "0xe"
$sheet2.Cells.Item(50, 6).Value = 380
$sheet2.Cells.Item(50, 7).Value = [double]"5.68432987514711e+23"
$sheet2.Cells.Item(50, 8).Value = 356
$sheet2.Cells.Item(50, 9).Value = 14

$sheet3 = $wb.Worksheets.Item("DE_PLT_#1")
$sheet3.Cells.Item(50, 1).Value = $rowDate
$sheet3.Cells.Item(50, 1).NumberFormat = $sheet3.Cells.Item(49, 1).NumberFormat
$sheet3.Cells.Item(50, 2).Value = "0x00,0x82"
$sheet3.Cells.Item(50, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$sheet3.Cells.Item(50, 4).Value = "0x00,0x7F"
$sheet3.Cells.Item(50, 5).Value = "0x7"
$sheet3.Cells.Item(50, 6).Value = 130
$sheet3.Cells.Item(50, 7).Value = [double]"5.68631262647114e+23"
$sheet3.Cells.Item(50, 8).Value = 127
$sheet3.Cells.Item(50, 9).Value = 7

$sheet4 = $wb.Worksheets.Item("DE_PLT_#2")
$sheet4.Cells.Item(50, 1).Value = $rowDate
$sheet4.Cells.Item(50, 1).NumberFormat = $sheet4.Cells.Item(49, 1).NumberFormat
$sheet4.Cells.Item(50, 2).Value = "0x00,0x82"
$sheet4.Cells.Item(50, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$sheet4.Cells.Item(50, 4).Value = "0x00,0x7E"
$sheet4.Cells.Item(50, 5).Value = "0x3"
$sheet4.Cells.Item(50, 6).Value = 130
$sheet4.Cells.Item(50, 7).Value = [double]"9.85046333984776e+23"
$sheet4.Cells.Item(50, 8).Value = 126
$sheet4.Cells.Item(50, 9).Value = 3

Write-Output "done"
